$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "C223217"
$ws.Range("D6").Value = "C2827321"
